# The workbook originally has three sheets, in this tab order:
#   1. "EXCEL SHEET"    -> Hello Excel / lecxE olleH / lecxE olleH_duplicated
#   2. "new_sheeeeeet1" -> empty
#   3. "Sheet"          -> Hello Excel
#
# Target state (two sheets, in this tab order):
#   1. "Sheet"       -> Hello Excel (unchanged)
#   2. "EXCEL SHEET" -> Hello world / dlrow olleH / dlrow olleH_duplicated

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Remove the extra, empty "new_sheeeeeet1" sheet entirely.
$wb.Worksheets.Item("new_sheeeeeet1").Delete()

# 2. Reorder remaining sheets so "Sheet" comes before "EXCEL SHEET".
$wb.Worksheets.Item("Sheet").Move($wb.Worksheets.Item(1))

# 3. Update the content of the "EXCEL SHEET" worksheet.
$dataSheet = $wb.Worksheets.Item("EXCEL SHEET")
$dataSheet.Range("A1").Value = "Hello world"
$dataSheet.Range("A2").Value = "dlrow olleH"
$dataSheet.Range("A3").Value = "dlrow olleH_duplicated"
